$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("A17").Value = 'Em lógica de programação, qual estrutura é usada para repetir um bloco de código várias vezes?'
$ws.Range("B17").Value = 'Condição'
$ws.Range("C17").Value = 'Função'
$ws.Range("D17").Value = 'Laço de repetição'
$ws.Range("E17").Value = 'Variável'
$ws.Range("F17").Value = 3

# Row 18
$ws.Range("A18").Value = 'O que significa o operador ''=='' em linguagens como Python, C e VisualG?'
$ws.Range("B18").Value = 'Atribuição'
$ws.Range("C18").Value = 'Comparação de igualdade'
$ws.Range("D18").Value = 'Concatenação'
$ws.Range("E18").Value = 'Incremento'
$ws.Range("F18").Value = 2

# Row 19
$ws.Range("A19").Value = 'Em Python, qual é o resultado de: print(2 * 3 ** 2)?'
$ws.Range("B19").Value = "'36"
$ws.Range("C19").Value = "'18"
$ws.Range("D19").Value = "'12"
$ws.Range("E19").Value = "'24"
$ws.Range("F19").Value = 2

# Row 20
$ws.Range("A20").Value = 'Qual é o tipo de dado usado para armazenar valores verdadeiros ou falsos?'
$ws.Range("B20").Value = 'String'
$ws.Range("C20").Value = 'Inteiro'
$ws.Range("D20").Value = 'Booleano'
$ws.Range("E20").Value = 'Caractere'
$ws.Range("F20").Value = 3

# Row 21
$ws.Range("A21").Value = 'Em C, qual comando finaliza uma instrução?'
$ws.Range("B21").Value = 'Dois pontos'
$ws.Range("C21").Value = 'Ponto e vírgula'
$ws.Range("D21").Value = 'Vírgula'
$ws.Range("E21").Value = 'Aspas'
$ws.Range("F21").Value = 2

# Row 22
$ws.Range("A22").Value = 'Em VisualG, qual estrutura é usada para tomada de decisão?'
$ws.Range("B22").Value = 'Enquanto'
$ws.Range("C22").Value = 'Para'
$ws.Range("D22").Value = 'Repita'
$ws.Range("E22").Value = 'Se...Então...Senão'
$ws.Range("F22").Value = 4

# Row 23
$ws.Range("A23").Value = 'Qual é a saída do código Python: print(len(''programacao''))?'
$ws.Range("B23").Value = "'11"
$ws.Range("C23").Value = "'10"
$ws.Range("D23").Value = "'9"
$ws.Range("E23").Value = "'12"
$ws.Range("F23").Value = 1

# Row 24
$ws.Range("A24").Value = 'Em lógica de programação, uma variável serve para:'
$ws.Range("B24").Value = 'Armazenar dados temporários'
$ws.Range("C24").Value = 'Criar gráficos'
$ws.Range("D24").Value = 'Apagar funções'
$ws.Range("E24").Value = 'Organizar pastas'
$ws.Range("F24").Value = 1

# Row 25
$ws.Range("A25").Value = 'Em C, qual destas é uma estrutura de repetição?'
$ws.Range("B25").Value = 'switch'
$ws.Range("C25").Value = 'if'
$ws.Range("D25").Value = 'for'
$ws.Range("E25").Value = 'typedef'
$ws.Range("F25").Value = 3

# Row 26
$ws.Range("A26").Value = 'Em Python, qual destes operadores representa ''OU lógico''?'
$ws.Range("B26").Value = '&&'
$ws.Range("C26").Value = '||'
$ws.Range("D26").Value = 'or'
$ws.Range("E26").Value = '&'
$ws.Range("F26").Value = 3

# Row 27
$ws.Range("A27").Value = 'O que significa ''indentação'' em Python?'
$ws.Range("B27").Value = 'Uso de ponto e vírgula'
$ws.Range("C27").Value = 'Espaços para definir blocos de código'
$ws.Range("D27").Value = 'Fechar chaves'
$ws.Range("E27").Value = 'Declarar variáveis'
$ws.Range("F27").Value = 2

# Row 28
$ws.Range("A28").Value = 'Em C, qual palavra-chave é usada para declarar uma constante?'
$ws.Range("B28").Value = 'const'
$ws.Range("C28").Value = 'static'
$ws.Range("D28").Value = 'define'
$ws.Range("E28").Value = 'final'
$ws.Range("F28").Value = 1

# Row 29
$ws.Range("A29").Value = 'No VisualG, qual comando exibe uma mensagem na tela?'
$ws.Range("B29").Value = 'escrever()'
$ws.Range("C29").Value = 'mostrar()'
$ws.Range("D29").Value = 'imprima()'
$ws.Range("E29").Value = 'escreva()'
$ws.Range("F29").Value = 4

# Row 30
$ws.Range("A30").Value = 'O que o comando ''break'' faz em estruturas de repetição?'
$ws.Range("B30").Value = 'Reinicia o laço'
$ws.Range("C30").Value = 'Finaliza o laço'
$ws.Range("D30").Value = 'Pausa temporariamente'
$ws.Range("E30").Value = 'Repete o bloco atual'
$ws.Range("F30").Value = 2

# Row 31
$ws.Range("A31").Value = 'Qual destas estruturas representa um laço ''para'' em Python?'
$ws.Range("B31").Value = 'for i = 1 até 10'
$ws.Range("C31").Value = 'for(i=0;i<10;i++)'
$ws.Range("D31").Value = 'for i in range(10):'
$ws.Range("E31").Value = 'loop 10 times'
$ws.Range("F31").Value = 3

# Row 32
$ws.Range("A32").Value = 'Qual é o resultado de: print(10 // 3) em Python?'
$ws.Range("B32").Value = "'3.3"
$ws.Range("C32").Value = "'3"
$ws.Range("D32").Value = "'4"
$ws.Range("E32").Value = "'3.0"
$ws.Range("F32").Value = 2

# Row 33
$ws.Range("A33").Value = 'Em C, qual operador é usado para acessar itens de um array?'
$ws.Range("B33").Value = '()'
$ws.Range("C33").Value = '[]'
$ws.Range("D33").Value = '{}'
$ws.Range("E33").Value = '<>'
$ws.Range("F33").Value = 2

# Row 34
$ws.Range("A34").Value = 'Em VisualG, qual destas é a estrutura correta do ''enquanto''?'
$ws.Range("B34").Value = 'enquanto(condição) { }'
$ws.Range("C34").Value = 'enquanto condição faça'
$ws.Range("D34").Value = 'while(condição)'
$ws.Range("E34").Value = 'loop enquanto'
$ws.Range("F34").Value = 2

# Row 35
$ws.Range("A35").Value = 'Em Python, qual é a função usada para ler entrada do usuário?'
$ws.Range("B35").Value = 'scan()'
$ws.Range("C35").Value = 'input()'
$ws.Range("D35").Value = 'read()'
$ws.Range("E35").Value = 'escreva()'
$ws.Range("F35").Value = 2

# Row 36
$ws.Range("A36").Value = 'Qual o valor final da variável x no código em C: int x=5; x+=3;?'
$ws.Range("B36").Value = "'2"
$ws.Range("C36").Value = "'8"
$ws.Range("D36").Value = "'15"
$ws.Range("E36").Value = "'3"
$ws.Range("F36").Value = 2

# Row 37
$ws.Range("A37").Value = 'No VisualG, qual comando inicia um algoritmo?'
$ws.Range("B37").Value = 'início'
$ws.Range("C37").Value = 'algoritmo'
$ws.Range("D37").Value = 'inicio'
$ws.Range("E37").Value = 'programa'
$ws.Range("F37").Value = 2

# Row 38
$ws.Range("A38").Value = 'Em Python, qual destas estruturas representa uma condição?'
$ws.Range("B38").Value = 'if x > 10:'
$ws.Range("C38").Value = 'se (x > 10)'
$ws.Range("D38").Value = 'if (x > 10)'
$ws.Range("E38").Value = 'condição(x > 10)'
$ws.Range("F38").Value = 1

# Row 39
$ws.Range("A39").Value = 'Em C, o que significa ''&&''?'
$ws.Range("B39").Value = 'OU lógico'
$ws.Range("C39").Value = 'NÃO lógico'
$ws.Range("D39").Value = 'E lógico'
$ws.Range("E39").Value = 'Comparação'
$ws.Range("F39").Value = 3

# Row 40
$ws.Range("A40").Value = 'Qual é a saída de: print(5 != 5) em Python?'
$ws.Range("B40").Value = "'True"
$ws.Range("C40").Value = "'False"
$ws.Range("D40").Value = 'Erro'
$ws.Range("E40").Value = "'5"
$ws.Range("F40").Value = 2

# Row 41
$ws.Range("A41").Value = 'Em lógica de programação, um algoritmo deve ser:'
$ws.Range("B41").Value = 'Confuso e longo'
$ws.Range("C41").Value = 'Ambíguo'
$ws.Range("D41").Value = 'Sequência de passos claros'
$ws.Range("E41").Value = 'Sempre recursivo'
$ws.Range("F41").Value = 3

Write-Output "Added rows 17-41"